$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'68.461.46"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +2.09%  "
$c = $ws.Range("D3")
$c.Value = "'3.649.97"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +1.42%  "
$ws.Range("E4").Value = "  -0.31%  "
$c = $ws.Range("D5")
$c.Value = "'195.80"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +9.28%  "
$c = $ws.Range("D6")
$c.Value = "'581.33"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.17%  "
$c = $ws.Range("D7")
$c.Value = "'3.644.76"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +1.44%  "
$ws.Range("E8").Value = "  +2.34%  "
$c = $ws.Range("D9")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("E10").Value = "  +2.24%  "
$c = $ws.Range("D11")
$c.Value = "'58.00"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +9.54%  "
$ws.Range("E12").Value = "  +7.39%  "
$ws.Range("E13").Value = "  +19.05%  "
$c = $ws.Range("D14")
$c.Value = "'10.22"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +3.65%  "
$c = $ws.Range("D15")
$c.Value = "'4.232.51"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.08%  "
$c = $ws.Range("D16")
$c.Value = "'3.652.75"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +1.11%  "
$ws.Range("E17").Value = "  +0.50%  "
$c = $ws.Range("D18")
$c.Value = "'12.64"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +4.17%  "
$c = $ws.Range("D19")
$c.Value = "'68.394.99"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +2.23%  "
$c = $ws.Range("D20")
$c.Value = "'18.66"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.42%  "
$ws.Range("E21").Value = "  +3.09%  "
$c = $ws.Range("D22")
$c.Value = "'404.32"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +3.85%  "
$c = $ws.Range("D23")
$c.Value = "'12.82"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +25.50%  "
$ws.Range("E24").Value = "  -0.15%  "
$c = $ws.Range("D25")
$c.Value = "'86.43"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.81%  "
$ws.Range("E26").Value = "  +4.87%  "
$c = $ws.Range("D27")
$c.Value = "'12.73"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +4.63%  "
$ws.Range("E28").Value = "  +8.67%  "
$ws.Range("E29").Value = "  +0.80%  "
$c = $ws.Range("D30")
$c.Value = "'8.18"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +22.79%  "
$c = $ws.Range("D31")
$c.Value = "'9.23"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +3.58%  "
$c = $ws.Range("D32")
$c.Value = "'31.82"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +2.91%  "
$c = $ws.Range("D33")
$c.Value = "'692.62"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +19.20%  "
$c = $ws.Range("D34")
$c.Value = "'12.31"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +4.45%  "
$ws.Range("E35").Value = "  +6.97%  "
$c = $ws.Range("D36")
$c.Value = "'65.07"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.08%  "
$c = $ws.Range("D37")
$c.Value = "'42.95"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +4.60%  "
$c = $ws.Range("D38")
$c.Value = "'0.424"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +13.94%  "
$ws.Range("E39").Value = "  +0.02%  "
$c = $ws.Range("D40")
$c.Value = "'0.0₃0799"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +9.84%  "
$ws.Range("E41").Value = "  +20.97%  "
$c = $ws.Range("D42")
$c.Value = "'3.14"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +14.41%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D43")
$c.Value = "'0.136"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +4.31%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$c = $ws.Range("D44")
$c.Value = "'3.215.70"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +19.99%  "
$c = $ws.Range("D45")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.28%  "
$c = $ws.Range("D46")
$c.Value = "'2.95"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +33.50%  "
$c = $ws.Range("D47")
$c.Value = "'0.0425"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +4.08%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$c = $ws.Range("D48")
$c.Value = "'3.18"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +4.32%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D49")
$c.Value = "'0.133"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +3.56%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$c = $ws.Range("D50")
$c.Value = "'8.92"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +9.64%  "
$c = $ws.Range("D51")
$c.Value = "'143.40"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +5.68%  "
